$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 78495561
$ws.Range("C2").Value = 78648161
$ws.Range("D2").Value = 78760072
$ws.Range("E2").Value = 78760516
$ws.Range("F2").Value = 78684291
$ws.Range("G2").Value = 78625520
$ws.Range("H2").Value = 78323212
$ws.Range("I2").Value = 77903290
$ws.Range("B3").Value = 21086389
$ws.Range("C3").Value = 21164659
$ws.Range("D3").Value = 21210964
$ws.Range("E3").Value = 21251482
$ws.Range("F3").Value = 21278736
$ws.Range("G3").Value = 21293740
$ws.Range("H3").Value = 21195793
$ws.Range("I3").Value = 21090443
$ws.Range("B4").Value = 7453631
$ws.Range("C4").Value = 7472040
$ws.Range("D4").Value = 7494748
$ws.Range("E4").Value = 7517343
$ws.Range("F4").Value = 7527269
$ws.Range("G4").Value = 7532155
$ws.Range("H4").Value = 7512462
$ws.Range("I4").Value = 7484393
$ws.Range("B5").Value = 7497393
$ws.Range("C5").Value = 7521297
$ws.Range("D6").Value = 8803658.890000001
$ws.Range("E6").Value = 8807193.890000001
$ws.Range("F6").Value = 8811002.890000001
$ws.Range("G6").Value = 8813496.92
$ws.Range("H6").Value = 8821282.860000003
$ws.Range("I6").Value = 8788866.910000002
$ws.Range("B7").Value = 5088576.94
$ws.Range("C7").Value = 5115453.900000001
$ws.Range("D7").Value = 5142714.920000001
$ws.Range("E7").Value = 5164543.910000001
$ws.Range("F7").Value = 5183915.930000001
$ws.Range("G7").Value = 5210900.890000002
$ws.Range("H7").Value = 5227572.94
$ws.Range("I7").Value = 5237195.900000001
$ws.Range("B8").Value = 16020490.93
$ws.Range("C8").Value = 15995199.91
$ws.Range("D8").Value = 15968206.9
$ws.Range("E8").Value = 15906640.94
$ws.Range("F8").Value = 15823769.92
$ws.Range("G8").Value = 15759019.92
$ws.Range("H8").Value = 15644154.94
$ws.Range("I8").Value = 15503611.88
$ws.Range("B9").Value = 6560049.900000001
$ws.Range("C9").Value = 6577240.910000001
$ws.Range("D9").Value = 6595066.900000001
$ws.Range("E9").Value = 6598630.910000001
$ws.Range("F9").Value = 6593432.900000001
$ws.Range("G9").Value = 6595448.930000002
$ws.Range("H9").Value = 6579055.870000002
$ws.Range("I9").Value = 6554249.910000001
$ws.Range("B10").Value = 10323880.88
$ws.Range("C10").Value = 10329931.95
$ws.Range("D10").Value = 10330159.9
$ws.Range("E10").Value = 10309136.9
$ws.Range("F11").Value = 9204632.870000001
$ws.Range("G11").Value = 9172853.91
$ws.Range("H11").Value = 9115555.91
$ws.Range("I11").Value = 9043290.9
$ws.Range("B12").Value = 3226018.91
$ws.Range("C12").Value = 3220047.91
$ws.Range("D12").Value = 3214535.920000001
$ws.Range("E12").Value = 3205527.890000001
$ws.Range("F13").Value = 4261519.930000001
$ws.Range("G13").Value = 4247891.880000002
$ws.Range("H13").Value = 4227324.91
$ws.Range("I13").Value = 4201221.91
$ws.Range("B14").Value = 1239110.93
$ws.Range("C14").Value = 1252271.91
